# Natmi following Dr Hou advice
# Rebuild the LR-pair table (rows 2-7) for every Sending-cluster x Target-cluster
# combination among {ECs, FAPs, sCs} x {FAPs, sCs}; Col1a2/Itga11 stay fixed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> FAPs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col1a2"
$ws.Range("C2").Value = "Itga11"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.256564333333333
$ws.Range("H2").Value = 15.769693
$ws.Range("I2").Value = 0.003747859920520347
$ws.Range("J2").Value = 0.003747859920520347
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 24.759128
$ws.Range("N2").Value = 74.277384
$ws.Range("O2").Value = 0.9895671066967037
$ws.Range("P2").Value = 0.9895671066967037
$ws.Range("Q2").Value = 130.1479491692347
$ws.Range("R2").Value = 1171.331542523112
$ws.Range("S2").Value = 0.003708758897853858
$ws.Range("T2").Value = 0.003708758897853857

# Row 3: ECs -> sCs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col1a2"
$ws.Range("C3").Value = "Itga11"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.256564333333333
$ws.Range("H3").Value = 15.769693
$ws.Range("I3").Value = 0.003747859920520347
$ws.Range("J3").Value = 0.003747859920520347
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.2610326666666667
$ws.Range("N3").Value = 0.7830980000000001
$ws.Range("O3").Value = 0.0104328933032964
$ws.Range("P3").Value = 0.0104328933032964
$ws.Range("Q3").Value = 1.372135005434889
$ws.Range("R3").Value = 12.349215048914
$ws.Range("S3").Value = 0.00003910102266648972
$ws.Range("T3").Value = 0.00003910102266648971

# Row 4: FAPs -> FAPs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Col1a2"
$ws.Range("C4").Value = "Itga11"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1312.703450666667
$ws.Range("H4").Value = 3938.110352
$ws.Range("I4").Value = 0.93593996730609
$ws.Range("J4").Value = 0.9359399673060897
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 24.759128
$ws.Range("N4").Value = 74.277384
$ws.Range("O4").Value = 0.9895671066967037
$ws.Range("P4").Value = 0.9895671066967037
$ws.Range("Q4").Value = 32501.39276109769
$ws.Range("R4").Value = 292512.5348498792
$ws.Range("S4").Value = 0.9261754054888949
$ws.Range("T4").Value = 0.9261754054888947

# Row 5: FAPs -> sCs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Col1a2"
$ws.Range("C5").Value = "Itga11"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1312.703450666667
$ws.Range("H5").Value = 3938.110352
$ws.Range("I5").Value = 0.93593996730609
$ws.Range("J5").Value = 0.9359399673060897
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.2610326666666667
$ws.Range("N5").Value = 0.7830980000000001
$ws.Range("O5").Value = 0.0104328933032964
$ws.Range("P5").Value = 0.0104328933032964
$ws.Range("Q5").Value = 342.6584822700551
$ws.Range("R5").Value = 3083.926340430496
$ws.Range("S5").Value = 0.00976456181719516
$ws.Range("T5").Value = 0.009764561817195157

# Row 6: sCs -> FAPs
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Col1a2"
$ws.Range("C6").Value = "Itga11"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 84.59089266666666
$ws.Range("H6").Value = 253.772678
$ws.Range("I6").Value = 0.06031217277338979
$ws.Range("J6").Value = 0.06031217277338978
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 24.759128
$ws.Range("N6").Value = 74.277384
$ws.Range("O6").Value = 0.9895671066967037
$ws.Range("P6").Value = 0.9895671066967037
$ws.Range("Q6").Value = 2094.396739168261
$ws.Range("R6").Value = 18849.57065251435
$ws.Range("S6").Value = 0.05968294230995504
$ws.Range("T6").Value = 0.05968294230995503

# Row 7: sCs -> sCs
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Col1a2"
$ws.Range("C7").Value = "Itga11"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 84.59089266666666
$ws.Range("H7").Value = 253.772678
$ws.Range("I7").Value = 0.06031217277338979
$ws.Range("J7").Value = 0.06031217277338978
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.2610326666666667
$ws.Range("N7").Value = 0.7830980000000001
$ws.Range("O7").Value = 0.0104328933032964
$ws.Range("P7").Value = 0.0104328933032964
$ws.Range("Q7").Value = 22.08098628849378
$ws.Range("R7").Value = 198.728876596444
$ws.Range("S7").Value = 0.0006292304634347539
$ws.Range("T7").Value = 0.0006292304634347538
